$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A and append a new entry right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$newRow = $lastRow + 1

# Carry over the formatting (borders, wrap text, etc.) from the row above
# so the new row matches the rest of the table's style.
$ws.Range("A" + $lastRow + ":B" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = "20-11-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹12,469 per gram for 24 karat gold, ₹11,430 per gram for 22 karat gold and ₹9,352 per gram for 18 karat gold (also called 999 gold)."
